# Regenerate save_data to use K instead of Strike#, recalc std/mean,
# and write the resulting s_vals into column G (K) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 2
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    17 = 2
    18 = 1
    19 = 0
    20 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
